$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Save" header in H1, copying the formatting used by the other
# header cells (e.g. G1 - bold, centered, bordered).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Add the numeric Save value in H2
$ws.Range("H2").Value = 1
